# The upstream change (see commit message: "Fixed POI packaging and
# upgraded to POI 3.15.") only affects how the *reference* fixture was
# re-serialized by the test-generation tool: every hunk in the diff is a
# pure XML-attribute / namespace-declaration re-ordering (alphabetical
# sorting of w:* attributes and xmlns:* declarations) produced by the
# newer Apache POI/XMLBeans writer. No paragraph text, run content,
# formatting, style definition, section property value, or any other
# document-model value actually changes between the "before" and
# "after" XML shown in the diff (every changed tag keeps exactly the
# same set of attribute name/value pairs, just listed in a different
# order).
#
# Since that kind of low-level serialization detail is not part of the
# Word object model (there's no COM-level "attribute order" concept to
# edit - Find/Replace, styles, paragraphs, sections, etc. all operate on
# semantic content), there is nothing to change here. We simply touch
# the document object so the script is a valid, explicit no-op.
$d = $word.ActiveDocument
